# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to Sheet1,
# populating headers and data for all existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used data row (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New header cells, matching the bold/centered header style used by the
# existing headers (copy style from L1, the last existing header cell).
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20141190
    $ws.Cells.Item($r, 15).Value = 10
}
